$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").Value = "'2"
$ws.Range("D2").Value = 0.12718
$ws.Range("F2").Value = 0.304
$ws.Range("G2").Value = 0.3440944881889763
$ws.Range("H2").Value = -0.3614173228346457
$ws.Range("I2").Value = -0.4072034700605307
$ws.Range("J2").Value = -0.4072034700605307
$ws.Range("K2").Value = -6.68
$ws.Range("L2").Value = -0.7514060742407198
$ws.Range("M2").Value = 0.135
$ws.Range("N2").Value = 0.003537735849056604
$ws.Range("O2").Value = -0.02020958083832336
$ws.Range("S2").Value = 0.135
$ws.Range("U2").Value = 3.602
$ws.Range("V2").Value = 0.09439203354297694
$ws.Range("W2").Value = -0.3484783601690076
$ws.Range("X2").Value = 0.09485093502015862
$ws.Range("Y2").Value = -0.4433292951891662
$ws.Range("Z2").Value = 0.3588121691483948
$ws.Range("AA2").Value = -0.1357866488304114
$ws.Range("AB2").Value = 0.07080491912525522
$ws.Range("AC2").Value = -0.2065915679556666
$ws.Range("AD2").Value = 8.92
$ws.Range("AE2").Value = 0.4601942441905907
$ws.Range("AF2").Value = 9.380194244190591
$ws.Range("AG2").Value = 5.77819424419059
$ws.Range("AH2").Value = 0.1973108102169113
$ws.Range("AI2").Value = 0.3526363062645379
$ws.Range("AJ2").Value = 0.131507321672751
$ws.Range("AK2").Value = 0.2512455622749677
$ws.Range("AL2").Value = 0.314
$ws.Range("AM2").Value = 0.314
$ws.Range("AN2").Value = -3.006403774856757
$ws.Range("AO2").Value = -12.11146496815287
$ws.Range("AP2").Value = -1.947487106232083
$ws.Range("AQ2").Value = -12.11146496815287

# --- Row 3 updates ---
$ws.Range("D3").Value = -0.00764
$ws.Range("G3").Value = -0.08892988929889299
$ws.Range("H3").Value = -0.274169741697417
$ws.Range("I3").Value = -0.3070989735960191
$ws.Range("J3").Value = -0.3070989735960191
$ws.Range("K3").Value = -1.25
$ws.Range("L3").Value = -0.4612546125461255
$ws.Range("M3").Value = 0.135
$ws.Range("N3").Value = 0.0266798418972332
$ws.Range("O3").Value = -0.108
$ws.Range("S3").Value = 0.135
$ws.Range("U3").Value = 3.22
$ws.Range("V3").Value = 0.6363636363636365
$ws.Range("W3").Value = -0.1388888888888889
$ws.Range("X3").Value = 0.1131695098436263
$ws.Range("Y3").Value = -0.2520583987325152
$ws.Range("Z3").Value = 0.2730150945937767
$ws.Range("AA3").Value = -0.08384265532596887
$ws.Range("AB3").Value = 0.0703949260840378
$ws.Range("AC3").Value = -0.1542375814100067
$ws.Range("AD3").Value = 4.74
$ws.Range("AE3").Value = 0.02619109222605863
$ws.Range("AF3").Value = 4.766191092226059
$ws.Range("AG3").Value = 1.546191092226059
$ws.Range("AH3").Value = 0.4850497051697689
$ws.Range("AI3").Value = 0.3765896909658424
$ws.Range("AJ3").Value = 0.2340518266335898
$ws.Range("AK3").Value = 0.1638575434848789
$ws.Range("AL3").Value = 0.1
$ws.Range("AM3").Value = 0.1
$ws.Range("AN3").Value = -6.26984126984127
$ws.Range("AO3").Value = -8.43
$ws.Range("AP3").Value = -2.045226312468332
$ws.Range("AQ3").Value = -8.43

# --- Row 4 (new row) ---
$ws.Range("A4").Value = "Spain"
$ws.Range("B4").Value = "Pangaea Oncology, S.A. (BME:PANG)"
$ws.Range("C4").Value = "Heathcare Information and Technology"
$ws.Range("D4").Value = 0.262
$ws.Range("F4").Value = 0.304
$ws.Range("G4").Value = 0.5339805825242717
$ws.Range("H4").Value = -0.3996763754045308
$ws.Range("I4").Value = -0.4511004256299202
$ws.Range("J4").Value = -0.4511004256299202
$ws.Range("K4").Value = -5.43
$ws.Range("L4").Value = -0.8786407766990291
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.382
$ws.Range("V4").Value = 0.01154078549848943
$ws.Range("W4").Value = -0.5580678314491263
$ws.Range("X4").Value = 0.07653236019669095
$ws.Range("Y4").Value = -0.6346001916458173
$ws.Range("Z4").Value = 0.4161615278298737
$ws.Range("AA4").Value = -0.1877306423348539
$ws.Range("AB4").Value = 0.07121491216647265
$ws.Range("AC4").Value = -0.2589455545013266
$ws.Range("AD4").Value = 4.18
$ws.Range("AE4").Value = 0.4340031519645321
$ws.Range("AF4").Value = 4.614003151964532
$ws.Range("AG4").Value = 4.232003151964532
$ws.Range("AH4").Value = 0.1223419092736695
$ws.Range("AI4").Value = 0.3308951598533222
$ws.Range("AJ4").Value = 0.1133612663305968
$ws.Range("AK4").Value = 0.312048530334658
$ws.Range("AL4").Value = 0.214
$ws.Range("AM4").Value = 0.214
$ws.Range("AN4").Value = -1.890547263681592
$ws.Range("AO4").Value = -13.83177570093458
$ws.Range("AP4").Value = -1.914067459052253
$ws.Range("AQ4").Value = -13.83177570093458
